# 2023 Day 20 done - update stats sheet and Overall tracker.

$wb = $excel.ActiveWorkbook

$ws2023 = $wb.Worksheets.Item("2023")
$wsOverall = $wb.Worksheets.Item("Overall")

# --- Update daily B/C (My 1 / My 2) values for days that shifted totals ---
$updates = @(
    @{ Row = 2;  B = 219499; C = 70141 },
    @{ Row = 3;  B = 185753; C = 8664 },
    @{ Row = 4;  B = 122981; C = 18530 },
    @{ Row = 5;  B = 122501; C = 16529 },
    @{ Row = 6;  B = 75067;  C = 29113 },
    @{ Row = 7;  B = 96400;  C = 1669 },
    @{ Row = 8;  B = 75956;  C = 6725 },
    @{ Row = 9;  B = 68395;  C = 13479 },
    @{ Row = 10; B = 70145;  C = 1055 },
    @{ Row = 11; B = 43890;  C = 15635 },
    @{ Row = 12; B = 52032;  C = 2133 },
    @{ Row = 13; B = 27163;  C = 13545 },
    @{ Row = 14; B = 33882;  C = 4738 },
    @{ Row = 15; B = 32010;  C = 6701 },
    @{ Row = 16; B = 35737;  C = 3719 },
    @{ Row = 17; B = 30215;  C = 915 },
    @{ Row = 18; B = 19003;  C = 1046 },
    @{ Row = 19; B = 19213;  C = 4641 },
    @{ Row = 20; B = 16216;  C = 6636 }
)

foreach ($u in $updates) {
    $ws2023.Cells.Item($u.Row, 2).Value = $u.B
    $ws2023.Cells.Item($u.Row, 3).Value = $u.C
}

# --- Day 20 (row 21) newly filled in ---
# (C21 is written before B21: the D21/"=B21+C21" formula's ISBLANK(B21) guard
#  only flips to "not blank" when B21 is written, so B21 must be the last of
#  the pair to land for the dependent formulas to pick up both new values.)
$ws2023.Cells.Item(21, 3).Value = 4077    # C21 My 2
$ws2023.Cells.Item(21, 5).Value = 9863    # E21 Overall 1
$ws2023.Cells.Item(21, 6).Value = 7499    # F21 Overall 2
$ws2023.Cells.Item(21, 2).Value = 8907    # B21 My 1

# H21 (=IF(ISBLANK(C21),"",F21/B21)) briefly divides by a zero/blank B21 while
# the values above land one at a time; re-stamping its formula forces the
# shared-formula group to recompute from the final settled inputs instead of
# keeping the transient #DIV/0! it cached mid-update.
$ws2023.Cells.Item(21, 8).Formula = '=IF(ISBLANK(C21),"",F21/B21)'

# --- Mark Day 20 as succeeded ("s") for 2023 on the Overall tracker ---
$wsOverall.Range("BZ12:CC12").Value = "s"

# --- Selections / active sheet reflect where the author ended up editing ---
[void]$ws2023.Range("F22").Select()
[void]$wsOverall.Select()
[void]$wsOverall.Range("CW12").Select()
